$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HU_FertilityByYear")

# Rename sheet HU_FertilityByYear -> PL_FertilityByYear
$ws.Name = "PL_FertilityByYear"

# Update fertility values for row 2 (B2:AJ2) and apply integer number format
$ws.Range("B2:AJ2").NumberFormat = "0"

$ws.Range("B2").Value = 44.333333333333336
$ws.Range("C2").Value = 44.333333333333336
$ws.Range("D2").Value = 43
$ws.Range("E2").Value = 44.000000000000007
$ws.Range("F2").Value = 44.000000000000007
$ws.Range("G2").Value = 46.333333333333329
$ws.Range("H2").Value = 49.333333333333336
$ws.Range("I2").Value = 48.666666666666664
$ws.Range("J2").Value = 48
$ws.Range("K2").Value = 46.333333333333329
$ws.Range("L2").Value = 44.333333333333336
$ws.Range("M2").Value = 43
$ws.Range("N2").Value = 46.510101010101003
$ws.Range("O2").Value = 46.669774669774696
$ws.Range("P2").Value = 47.205000000000005
$ws.Range("Q2").Value = 47.089269101269103
$ws.Range("R2").Value = 47.258331552706601
$ws.Range("S2").Value = 47.427394004143999
$ws.Range("T2").Value = 47.596456455581503
$ws.Range("U2").Value = 48.467999999999996
$ws.Range("V2").Value = 48.075077577052603
$ws.Range("W2").Value = 48.254175472675499
$ws.Range("X2").Value = 48.433273368298401
$ws.Range("Y2").Value = 48.612371263921297
$ws.Range("Z2").Value = 49.54
$ws.Range("AA2").Value = 49.090331989639999
$ws.Range("AB2").Value = 49.276339400713297
$ws.Range("AC2").Value = 49.462346811786503
$ws.Range("AD2").Value = 49.648354222859801
$ws.Range("AE2").Value = 50.450333333333333
$ws.Range("AF2").Value = 50.102498604926303
$ws.Range("AG2").Value = 50.292480026963503
$ws.Range("AH2").Value = 50.482461449000603
$ws.Range("AI2").Value = 50.672442871037703
$ws.Range("AJ2").Value = 50.862424293074902
